$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "neni hotovy" text note that used to live in I2
$ws.Range("I2").ClearContents()

# Update column H hour values
$ws.Range("H2").Value = 9
$ws.Range("H10").Value = 6
$ws.Range("H14").Value = 7
$ws.Range("H17").Value = 10
$ws.Range("H20").Value = 3

# Add total formula
$ws.Range("H22").Formula = "=SUM(H2:H21)"

# Update selection to reflect the new active cell used while editing
$ws.Range("I22").Select()
